$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows before row 16. This pushes the existing footer rows
# (old A21 "NSB" / old A22 citation) down to rows 27/28, leaving rows
# 16-20 free for the new employee-count / assets / turnover breakdown
# table and rows 21-26 blank - matching the target layout where the
# footer now lives at A27/A28 instead of A21/A22.
$ws.Range("A16:A21").EntireRow.Insert()

# Row 16: new table header ("Number of employees" / "Assets ..." /
# "Turnover ..."), bold - same emphasis as the other section headers
# already in the sheet (e.g. A3, A7, A10-A12).
$ws.Range("B16").Value = "Number of employees"
$ws.Range("C16").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D16").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B16:D16").Font.Bold = $true

# Rows 17-20: the Micro/Small/Medium/Large employee-count breakdown.
# Columns C/D are intentionally left blank (no Assets/Turnover figures
# reported for Bhutan), matching the source data.
$ws.Range("A17").Value = "Micro"
$ws.Range("B17").Value = "1-4"

$ws.Range("A18").Value = "Small"
$ws.Range("B18").Value = "5-19"

$ws.Range("A19").Value = "Medium"
$ws.Range("B19").Value = "20-99"

$ws.Range("A20").Value = "Large"
$ws.Range("B20").Value = ">=100"

# The row-insert carries the old A21/A22 formatting along for the ride,
# but re-assert it explicitly on its new home (A27/A28) so the "NSB"
# source label stays bold and the citation stays italic, same as before
# the rows moved.
$ws.Range("A27").Font.Bold = $true
$ws.Range("A28").Font.Italic = $true
